# edit.ps1
# Applies the "updated high level TOC" change to the brief-contents TOC document.
#
# Summary of changes:
#  1. "Part 1:" line gains a new highlighted subtitle run: " Hypermedia: The New/Old Way"
#  2. "Part 2:" line gains a new highlighted subtitle run: " Modern Hypermedia"
#  3. "Part 3:" line gains a new highlighted subtitle run: " Hypermedia & Other Technologies"
#  4. A new "Part 4: Developing With Hypermedia " heading paragraph is inserted
#     immediately before the "10 Creating A Dynamic Download UI" chapter line
#     (gaining a trailing space on that chapter line), and the old, bare
#     "Part 4: " heading paragraph (which used to sit after that chapter line) is removed.
#  5. The "13 Hypermedia: The New/Old Way " chapter line is simplified to
#     "13 Hypermedia Reconsidered".

$d = $word.ActiveDocument

function Insert-RunXml($range, $runsXml) {
    # Replaces the contents of $range (a non-collapsed Range) with one or more
    # freshly-formatted <w:r> runs, supplied as a raw OOXML fragment string.
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $runsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

function Insert-ParaXml($range, $pPrXml, $runsXml) {
    # Same idea as Insert-RunXml but also supplies paragraph properties (pPr).
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p>' + $pPrXml + $runsXml + '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $range.InsertXML($pkg)
}

# ---------------------------------------------------------------------------
# 1. "Part 1:" -> "Part 1: Hypermedia: The New/Old Way"
#    Replace the single trailing space run after "Part 1:" with a highlighted
#    run that carries the new subtitle text.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("Part 1: ") | Out-Null
$spaceRng = $d.Range($rng.End - 1, $rng.End)
$xml = '<w:r><w:rPr><w:shd w:val="clear" w:color="auto" w:fill="ffff00"/><w:rtl w:val="0"/><w:lang w:val="en-US"/></w:rPr>' +
       '<w:t xml:space="preserve"> Hypermedia: The New/Old Way</w:t></w:r>'
Insert-RunXml $spaceRng $xml

# ---------------------------------------------------------------------------
# 2. "Part 2:" -> "Part 2: Modern Hypermedia"
#    Replace the trailing space run with [new highlighted run][original space run].
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("Part 2: ") | Out-Null
$spaceRng = $d.Range($rng.End - 1, $rng.End)
$xml = '<w:r><w:rPr><w:shd w:val="clear" w:color="auto" w:fill="ffff00"/><w:rtl w:val="0"/><w:lang w:val="en-US"/></w:rPr>' +
       '<w:t xml:space="preserve"> Modern Hypermedia</w:t></w:r>' +
       '<w:r><w:rPr><w:rtl w:val="0"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>'
Insert-RunXml $spaceRng $xml

# ---------------------------------------------------------------------------
# 3. "Part 3:" -> "Part 3: Hypermedia & Other Technologies"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("Part 3: ") | Out-Null
$spaceRng = $d.Range($rng.End - 1, $rng.End)
$xml = '<w:r><w:rPr><w:shd w:val="clear" w:color="auto" w:fill="ffff00"/><w:rtl w:val="0"/><w:lang w:val="en-US"/></w:rPr>' +
       '<w:t xml:space="preserve"> Hypermedia &amp; Other Technologies</w:t></w:r>' +
       '<w:r><w:rPr><w:rtl w:val="0"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>'
Insert-RunXml $spaceRng $xml

# ---------------------------------------------------------------------------
# 4. Move/rewrite the "Part 4:" heading:
#    a) add a trailing space to the "10 Creating A Dynamic Download UI" line
#    b) delete the old bare "Part 4: " paragraph that used to follow it
#    c) insert a new, highlighted "Part 4: Developing With Hypermedia " heading
#       paragraph immediately before that chapter line
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("10 Creating A Dynamic Download UI") | Out-Null
$chapterStart = $rng.Start
$chapterEnd = $rng.End

$endPos = $d.Range($chapterEnd, $chapterEnd)
$endPos.InsertAfter(" ")

$afterChapter = $d.Range($chapterEnd, $d.Content.End)
$afterChapter.Find.ClearFormatting()
$afterChapter.Find.Execute("Part 4:") | Out-Null
$oldPart4Para = $afterChapter.Paragraphs(1)
$oldPart4Para.Range.Delete()

$insPos = $d.Range($chapterStart, $chapterStart)
$insPos.InsertParagraphBefore()
$newParaRange = $d.Range($chapterStart, $chapterStart + 1)
$pPr = '<w:pPr><w:pStyle w:val="TOC parts"/></w:pPr>'
$runs = '<w:r><w:rPr><w:shd w:val="clear" w:color="auto" w:fill="ffff00"/><w:rtl w:val="0"/><w:lang w:val="en-US"/></w:rPr><w:t>Part 4:</w:t></w:r>' +
        '<w:r><w:rPr><w:shd w:val="clear" w:color="auto" w:fill="ffff00"/><w:rtl w:val="0"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Developing With </w:t></w:r>' +
        '<w:r><w:rPr><w:shd w:val="clear" w:color="auto" w:fill="ffff00"/><w:rtl w:val="0"/><w:lang w:val="en-US"/></w:rPr><w:t>Hypermedia</w:t></w:r>' +
        '<w:r><w:rPr><w:rtl w:val="0"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>'
Insert-ParaXml $newParaRange $pPr $runs

# ---------------------------------------------------------------------------
# 5. "13 Hypermedia: The New/Old Way " -> "13 Hypermedia Reconsidered"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute("13 Hypermedia: ") | Out-Null
$paraRng = $rng.Paragraphs(1)
$contentRng = $d.Range($paraRng.Range.Start, $paraRng.Range.End - 1)
$xml = '<w:r><w:rPr><w:rtl w:val="0"/><w:lang w:val="it-IT"/></w:rPr><w:t>13 Hypermedia</w:t></w:r>' +
       '<w:r><w:rPr><w:rtl w:val="0"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> Reconsidered</w:t></w:r>'
Insert-RunXml $contentRng $xml
